# Update "想去人数" (F column) values on the sheets that list exhibition
# events: "展览" (sheet1) and "全部类型" (sheet4). The two sheets share the
# same set of events, but sheet4 has one extra row inserted partway through,
# so the row numbers for the later entries differ between the two sheets.

$wb = $excel.ActiveWorkbook

$sheetUpdates = @{
    "展览"    = @{ "F5" = 3058; "F6" = 2043; "F7" = 396; "F9"  = 1139; "F11" = 825; "F12" = 72 }
    "全部类型" = @{ "F5" = 3058; "F6" = 2043; "F7" = 396; "F10" = 1139; "F12" = 825; "F13" = 72 }
}

foreach ($sheetName in $sheetUpdates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $updates = $sheetUpdates[$sheetName]
    foreach ($cellRef in $updates.Keys) {
        $ws.Range($cellRef).Value = $updates[$cellRef]
    }
}
